# issue #5: stock data from json to db
# The "股票" (stock) worksheet gains three new columns: a "category" column
# inserted right after "property_category", and "source_file"/"index"
# columns appended at the end. The existing data row is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "face_value"
$ws.Range("F1").Value = "currency"
$ws.Range("G1").Value = "total"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2) ---
$ws.Range("A2").Value = 73
$ws.Range("B2").Value = "馬祖酒廠實業股份有限公司(未上市）"
$ws.Range("C2").Value = "陳雪生"
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "新臺幣"
$ws.Range("G2").Value = 10000
$ws.Range("H2").Value = "stock"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-27"
$ws.Range("K2").Value = "陳雪生"
$ws.Range("L2").Value = 1751
$ws.Range("M2").Value = "tmp5a001"
$ws.Range("N2").Value = 73

# --- Copy formatting for the newly-introduced cells so they match the
# --- look of the rest of their row (header style vs. data style). ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("L2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
